$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-28)
# from 45208 to 45212
$ws.Range("C2:C28").Value = 45212

# Update row 2 hyperlink formulas to point to the renamed files
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2482/artfynd/A 30703-2023 artfynd.xlsx", "A 30703-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2482/kartor/A 30703-2023 karta.png", "A 30703-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2482/klagomål/A 30703-2023 fsc-klagomål.docx", "A 30703-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2482/klagomålsmail/A 30703-2023 fsc-klagomål mail.docx", "A 30703-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2482/tillsyn/A 30703-2023 tillsynsbegäran.docx", "A 30703-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2482/ti,llsynsmail/A 30703-2023 tillsynsbegäran mail.docx", "A 30703-2023")'
